# Apply "Remove unnecessary activities and variables" changes to Config.xlsx
#
# Summary of content changes (derived from the OOXML diff):
#  - Settings sheet: remove the OrchestratorQueueName / OrchestratorQueueFolder /
#    logF_BusinessProcessName rows (rows 2, 3 and 5 become empty, keeping their
#    existing cell styles).
#  - Constants sheet: MaxRetryNumber value 0 -> 2, MaxConsecutiveSystemExceptions
#    value 0 -> 3, ShouldMarkJobAsFaulted value FALSE -> TRUE.
#  - View/selection bookkeeping: Settings becomes the active/selected sheet with
#    B16 selected, Constants keeps C7 selected, Assets loses the "tabSelected" flag.

$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# --- Settings sheet: drop the queue-related and business-process-name rows ---
# Row 2: OrchestratorQueueName | ProcessABCQueue | <description>
$wsSettings.Range("A2:C2").ClearContents()
# Row 3: OrchestratorQueueFolder | <empty> | <description>
$wsSettings.Range("A3:C3").ClearContents()
# Row 5: logF_BusinessProcessName | Framework | <description>
$wsSettings.Range("A5:C5").ClearContents()

# The removed text no longer needs to wrap across multiple lines, so the rows
# shrink back down to a single text line once their content is gone.
$wsSettings.Rows.Item(3).RowHeight = 14.4
$wsSettings.Rows.Item(5).RowHeight = 14.4

# --- Constants sheet: updated default values ---
$wsConstants.Range("B2").Value = 2
$wsConstants.Range("B3").Value = 3
$wsConstants.Range("B17").Value = $true

# --- Selection / active sheet bookkeeping ---
$wsConstants.Activate()
$wsConstants.Range("C7").Select()

$wsSettings.Activate()
$wsSettings.Range("B16").Select()
